# Revert "Drop in all data files from 3.0 RMI script"
# Change the currency-year references from 2019 back to 2018 on the
# "About" sheet, and update the derived conversion-factor value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Text cells that mention the dollar "vintage" year: 2019 -> 2018
$ws.Range("A18").Value = "billion 2018 dollars"
$ws.Range("A21").Value = "million 2018 dollars"
$ws.Range("B26").Value = "2018 dollars per 2012 dollar"
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2018 dollar."'

# Updated conversion factor (2012 dollars per 2018 dollar)
$ws.Range("A26").Value = 0.9143273584567535

$wb.Application.Calculate()
